$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.198.35"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").Value = "1.785.15"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'225.73"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'31.81"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").Value = "'0.291"
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("D12").Value = "2.042.86"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").Value = "1.790.71"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").Value = "34.133.86"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("D18").Value = "'68.15"
$ws.Range("E18").Value = "  +2.73%  "
$ws.Range("D19").Value = "'246.35"
$ws.Range("E19").Value = "  +4.04%  "
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'10.95"
$ws.Range("E21").Value = "  +4.31%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").Value = "  +2.79%  "
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").Value = "'161.44"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("E26").Value = "  +3.01%  "
$ws.Range("D27").Value = "'16.32"
$ws.Range("E27").Value = "  +1.78%  "
$ws.Range("E28").Value = "  +2.23%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "'1.24"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("D31").Value = "'0.0519"
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("E33").Value = "  +4.43%  "
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("D35").Value = "1.448.16"
$ws.Range("E35").Value = "  +5.38%  "
$ws.Range("E36").Value = "  +1.70%  "
$ws.Range("E37").Value = "  +9.46%  "
$ws.Range("E38").Value = "  +4.33%  "
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("D40").Value = "'80.37"
$ws.Range("E40").Value = "  +3.81%  "
$ws.Range("D41").Value = "'2.38"
$ws.Range("E42").Value = "  +2.83%  "
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("D44").Value = "'13.50"
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("D45").Value = "'6.06"
$ws.Range("D46").Value = "'0.0507"
$ws.Range("E46").Value = "  +2.15%  "
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("D49").Value = "1.944.78"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("D50").Value = "'105.93"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("E51").Value = "  +0.19%  "
